$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C6").Value = 8623
$ws.Range("C7:C14").Value = 8583
$ws.Range("C15:C23").Value = 8444
$ws.Range("C24:C29").Value = 8430
$ws.Range("C30:C39").Value = 8145
$ws.Range("C40:C45").Value = 8111
$ws.Range("C46:C56").Value = 8087
$ws.Range("C57:C57").Value = 7930
$ws.Range("C58:C60").Value = 7691
$ws.Range("C61:C252").Value = 7622
